# Auto-generated edit script: updates cryptos price/volume table
# to match the "Updated cryptos list" GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.668.05"
$ws.Range("E2").Value = "  -1.96%  "

$ws.Range("D3").Value = "1.797.83"
$ws.Range("E3").Value = "  -1.70%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "231.88"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.16%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5881"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.08%  "

$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2770"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.45%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06790"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.36%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.26"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.67%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07523"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.64%  "

$ws.Range("D12").Value = "1.792.48"
$ws.Range("E12").Value = "  -2.14%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.805"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.42%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6192"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.04%  "

$ws.Range("D15").Value = "2.042.06"
$ws.Range("E15").Value = "  -1.63%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000009113"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -8.31%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "75.48"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.31%  "

$ws.Range("D18").Value = "28.653.16"
$ws.Range("E18").Value = "  -1.93%  "

$ws.Range("E19").Value = "  -5.96%  "

$ws.Range("E20").Value = "  +0.02%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "210.64"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.71%  "

$ws.Range("E22").Value = "  -1.22%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.837"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.27%  "

$ws.Range("E24").Value = "  +0.05%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.74"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.02%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.980"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.26%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1264"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.28%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.46"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.36%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.425"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.70%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.06112"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.59%  "

$ws.Range("E31").Value = "  -1.26%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.807"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.35%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.793"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.06%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.739"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.15%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.050"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.31%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6432"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.24%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.501"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.73%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.711"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.98%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.507"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.14%  "

$ws.Range("E40").Value = "  -1.67%  "

$ws.Range("D41").Value = "1.143.14"
$ws.Range("E41").Value = "  -6.55%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8837"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.53%  "

$ws.Range("E43").Value = "  +0.33%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.19"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.11%  "

$ws.Range("D45").Value = "1.951.26"
$ws.Range("E45").Value = "  -1.54%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "60.26"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.38%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000113"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.78%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.597"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.39%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.357"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.68%  "

$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05476"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.39%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4480"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.72%  "
